$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (price, volume label, hour) from the
# Thu Dec 15 13:12:30 UTC 2022 GitHub Actions data refresh.
$updates = @(
    @{ Cell = 'D2'; Value = '264.71' }
    @{ Cell = 'G2'; Value = '13' }
    @{ Cell = 'D3'; Value = '22.88' }
    @{ Cell = 'G3'; Value = '13' }
    @{ Cell = 'D4'; Value = '6.236' }
    @{ Cell = 'G4'; Value = '13' }
    @{ Cell = 'D5'; Value = '0.06169' }
    @{ Cell = 'G5'; Value = '13' }
    @{ Cell = 'D6'; Value = '3.558' }
    @{ Cell = 'G6'; Value = '13' }
    @{ Cell = 'D7'; Value = '6.742' }
    @{ Cell = 'G7'; Value = '13' }
    @{ Cell = 'D8'; Value = '1.357' }
    @{ Cell = 'G8'; Value = '13' }
    @{ Cell = 'D9'; Value = '0.8124' }
    @{ Cell = 'G9'; Value = '13' }
    @{ Cell = 'D10'; Value = '0.1598' }
    @{ Cell = 'G10'; Value = '13' }
    @{ Cell = 'D11'; Value = '0.08228' }
    @{ Cell = 'G11'; Value = '13' }
    @{ Cell = 'D12'; Value = '0.03387' }
    @{ Cell = 'G12'; Value = '13' }
    @{ Cell = 'D13'; Value = '0.03171' }
    @{ Cell = 'G13'; Value = '13' }
    @{ Cell = 'D14'; Value = '0.09247' }
    @{ Cell = 'G14'; Value = '13' }
    @{ Cell = 'D15'; Value = '3.911' }
    @{ Cell = 'G15'; Value = '13' }
    @{ Cell = 'G16'; Value = '13' }
    @{ Cell = 'D17'; Value = '0.04854' }
    @{ Cell = 'G17'; Value = '13' }
    @{ Cell = 'D18'; Value = '0.0006315' }
    @{ Cell = 'E18'; Value = '17OneONEWorstin24h' }
    @{ Cell = 'G18'; Value = '13' }
    @{ Cell = 'D19'; Value = '0.006234' }
    @{ Cell = 'G19'; Value = '13' }
    @{ Cell = 'D20'; Value = '0.001099' }
    @{ Cell = 'G20'; Value = '13' }
    @{ Cell = 'D21'; Value = '0.003224' }
    @{ Cell = 'G21'; Value = '13' }
    @{ Cell = 'G22'; Value = '13' }
    @{ Cell = 'D23'; Value = '3.698' }
    @{ Cell = 'G23'; Value = '13' }
    @{ Cell = 'D24'; Value = '2.265' }
    @{ Cell = 'G24'; Value = '13' }
    @{ Cell = 'G25'; Value = '13' }
    @{ Cell = 'G26'; Value = '13' }
    @{ Cell = 'D27'; Value = '0.0002683' }
    @{ Cell = 'G27'; Value = '13' }
    @{ Cell = 'G28'; Value = '13' }
    @{ Cell = 'G29'; Value = '13' }
    @{ Cell = 'G30'; Value = '13' }
    @{ Cell = 'G31'; Value = '13' }
    @{ Cell = 'G32'; Value = '13' }
    @{ Cell = 'G33'; Value = '13' }
    @{ Cell = 'G34'; Value = '13' }
    @{ Cell = 'G35'; Value = '13' }
    @{ Cell = 'G36'; Value = '13' }
    @{ Cell = 'G37'; Value = '13' }
    @{ Cell = 'G38'; Value = '13' }
    @{ Cell = 'G39'; Value = '13' }
    @{ Cell = 'D40'; Value = '0.04595' }
    @{ Cell = 'G40'; Value = '13' }
    @{ Cell = 'D41'; Value = '0.007401' }
    @{ Cell = 'G41'; Value = '13' }
    @{ Cell = 'G42'; Value = '13' }
    @{ Cell = 'D43'; Value = '0.003135' }
    @{ Cell = 'G43'; Value = '13' }
    @{ Cell = 'D44'; Value = '0.01081' }
    @{ Cell = 'E44'; Value = '43LocalTradersLCT' }
    @{ Cell = 'G44'; Value = '13' }
    @{ Cell = 'D45'; Value = '0.00006159' }
    @{ Cell = 'G45'; Value = '13' }
    @{ Cell = 'D46'; Value = '0.00000000750' }
    @{ Cell = 'G46'; Value = '13' }
    @{ Cell = 'D47'; Value = '0.7505' }
    @{ Cell = 'G47'; Value = '13' }
    @{ Cell = 'D48'; Value = '0.2523' }
    @{ Cell = 'G48'; Value = '13' }
    @{ Cell = 'D49'; Value = '0.00002101' }
    @{ Cell = 'G49'; Value = '13' }
    @{ Cell = 'G50'; Value = '13' }
    @{ Cell = 'G51'; Value = '13' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
